# chore(runtime): publish files + archive (2025-12-05 07:28:44)
# Applies the KHL stats refresh: appends 3 newly-completed matches to
# Matches_SOG, rolls the Shots_HA / Shots_Summary "as_of" snapshot forward
# to the 2025-12-04 16:30 UTC slate (refreshing the six rows for the teams
# that played), and bumps Meta_ext's as_of_utc / build_version.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Matches_SOG: append the 3 new match rows (354-356)
# ---------------------------------------------------------------------
$wsMatches = $wb.Worksheets.Item("Matches_SOG")

# Column A (uid) holds numeric-looking ids but is stored as TEXT in this
# sheet (matches the existing 2..353 rows) - force text formatting before
# assigning so the COM layer doesn't silently coerce it to a number.
$uidRange = $wsMatches.Range("A354:A356")
$uidRange.NumberFormat = "@"

$newMatches = @(
    @(354, "897848", "2025-12-04T12:15:00", "Амур",     "Локомотив", 21, 27),
    @(355, "897849", "2025-12-04T12:30:00", "Адмирал",  "СКА",       33, 33),
    @(356, "897847", "2025-12-04T16:30:00", "Авангард", "ХК Сочи",   31, 23)
)

foreach ($row in $newMatches) {
    $r = $row[0]
    $wsMatches.Cells.Item($r, 1).Value = $row[1]
    $wsMatches.Cells.Item($r, 2).Value = $row[2]
    $wsMatches.Cells.Item($r, 3).Value = $row[3]
    $wsMatches.Cells.Item($r, 4).Value = $row[4]
    $wsMatches.Cells.Item($r, 5).Value = $row[5]
    $wsMatches.Cells.Item($r, 6).Value = $row[6]
    $wsMatches.Cells.Item($r, 7).Value = "khl_text"
}

# ---------------------------------------------------------------------
# Shots_HA: roll as_of_utc forward for every team, refresh the six rows
# (home team + away team of each new match) with their updated totals.
# ---------------------------------------------------------------------
$wsHA = $wb.Worksheets.Item("Shots_HA")
$wsHA.Cells.Item(2, 4).Value = "2025-12-04T16:30:00Z"
$wsHA.Cells.Item(2, 5).Value = 15
$wsHA.Cells.Item(2, 7).Value = 480
$wsHA.Cells.Item(2, 8).Value = 427
$wsHA.Cells.Item(2, 9).Value = 32
$wsHA.Cells.Item(2, 10).Value = 28.5
$wsHA.Cells.Item(3, 4).Value = "2025-12-04T16:30:00Z"
$wsHA.Cells.Item(4, 4).Value = "2025-12-04T16:30:00Z"
$wsHA.Cells.Item(4, 5).Value = 15
$wsHA.Cells.Item(4, 7).Value = 555
$wsHA.Cells.Item(4, 8).Value = 408
$wsHA.Cells.Item(4, 9).Value = 37
$wsHA.Cells.Item(4, 10).Value = 27.2
$wsHA.Cells.Item(5, 4).Value = "2025-12-04T16:30:00Z"
$wsHA.Cells.Item(6, 4).Value = "2025-12-04T16:30:00Z"
$wsHA.Cells.Item(6, 5).Value = 17
$wsHA.Cells.Item(6, 7).Value = 506
$wsHA.Cells.Item(6, 8).Value = 581
$wsHA.Cells.Item(6, 9).Value = 29.8
$wsHA.Cells.Item(6, 10).Value = 34.2
$wsHA.Cells.Item(7, 4).Value = "2025-12-04T16:30:00Z"
$wsHA.Cells.Item(8, 4).Value = "2025-12-04T16:30:00Z"
$wsHA.Cells.Item(9, 4).Value = "2025-12-04T16:30:00Z"
$wsHA.Cells.Item(10, 4).Value = "2025-12-04T16:30:00Z"
$wsHA.Cells.Item(11, 4).Value = "2025-12-04T16:30:00Z"
$wsHA.Cells.Item(12, 4).Value = "2025-12-04T16:30:00Z"
$wsHA.Cells.Item(12, 6).Value = 18
$wsHA.Cells.Item(12, 11).Value = 551
$wsHA.Cells.Item(12, 12).Value = 441
$wsHA.Cells.Item(12, 13).Value = 30.6
$wsHA.Cells.Item(12, 14).Value = 24.5
$wsHA.Cells.Item(13, 4).Value = "2025-12-04T16:30:00Z"
$wsHA.Cells.Item(14, 4).Value = "2025-12-04T16:30:00Z"
$wsHA.Cells.Item(15, 4).Value = "2025-12-04T16:30:00Z"
$wsHA.Cells.Item(15, 6).Value = 16
$wsHA.Cells.Item(15, 11).Value = 476
$wsHA.Cells.Item(15, 12).Value = 541
$wsHA.Cells.Item(15, 13).Value = 29.8
$wsHA.Cells.Item(15, 14).Value = 33.8
$wsHA.Cells.Item(16, 4).Value = "2025-12-04T16:30:00Z"
$wsHA.Cells.Item(17, 4).Value = "2025-12-04T16:30:00Z"
$wsHA.Cells.Item(18, 4).Value = "2025-12-04T16:30:00Z"
$wsHA.Cells.Item(19, 4).Value = "2025-12-04T16:30:00Z"
$wsHA.Cells.Item(20, 4).Value = "2025-12-04T16:30:00Z"
$wsHA.Cells.Item(21, 4).Value = "2025-12-04T16:30:00Z"
$wsHA.Cells.Item(22, 4).Value = "2025-12-04T16:30:00Z"
$wsHA.Cells.Item(22, 6).Value = 16
$wsHA.Cells.Item(22, 11).Value = 415
$wsHA.Cells.Item(22, 12).Value = 588
$wsHA.Cells.Item(22, 13).Value = 25.9
$wsHA.Cells.Item(22, 14).Value = 36.8
$wsHA.Cells.Item(23, 4).Value = "2025-12-04T16:30:00Z"

# ---------------------------------------------------------------------
# Shots_Summary: same as_of_utc roll-forward + refreshed SOG totals for
# the six affected teams.
# ---------------------------------------------------------------------
$wsSum = $wb.Worksheets.Item("Shots_Summary")
$wsSum.Cells.Item(2, 4).Value = "2025-12-04T16:30:00Z"
$wsSum.Cells.Item(2, 5).Value = 30
$wsSum.Cells.Item(2, 6).Value = 1011
$wsSum.Cells.Item(2, 7).Value = 933
$wsSum.Cells.Item(2, 8).Value = 33.7
$wsSum.Cells.Item(2, 9).Value = 31.1
$wsSum.Cells.Item(3, 4).Value = "2025-12-04T16:30:00Z"
$wsSum.Cells.Item(4, 4).Value = "2025-12-04T16:30:00Z"
$wsSum.Cells.Item(4, 5).Value = 31
$wsSum.Cells.Item(4, 6).Value = 1052
$wsSum.Cells.Item(4, 7).Value = 852
$wsSum.Cells.Item(4, 8).Value = 33.9
$wsSum.Cells.Item(4, 9).Value = 27.5
$wsSum.Cells.Item(5, 4).Value = "2025-12-04T16:30:00Z"
$wsSum.Cells.Item(6, 4).Value = "2025-12-04T16:30:00Z"
$wsSum.Cells.Item(6, 5).Value = 33
$wsSum.Cells.Item(6, 6).Value = 944
$wsSum.Cells.Item(6, 7).Value = 1195
$wsSum.Cells.Item(6, 8).Value = 28.6
$wsSum.Cells.Item(6, 9).Value = 36.2
$wsSum.Cells.Item(7, 4).Value = "2025-12-04T16:30:00Z"
$wsSum.Cells.Item(8, 4).Value = "2025-12-04T16:30:00Z"
$wsSum.Cells.Item(9, 4).Value = "2025-12-04T16:30:00Z"
$wsSum.Cells.Item(10, 4).Value = "2025-12-04T16:30:00Z"
$wsSum.Cells.Item(11, 4).Value = "2025-12-04T16:30:00Z"
$wsSum.Cells.Item(12, 4).Value = "2025-12-04T16:30:00Z"
$wsSum.Cells.Item(12, 5).Value = 35
$wsSum.Cells.Item(12, 6).Value = 1123
$wsSum.Cells.Item(12, 7).Value = 891
$wsSum.Cells.Item(12, 8).Value = 32.1
$wsSum.Cells.Item(12, 9).Value = 25.5
$wsSum.Cells.Item(13, 4).Value = "2025-12-04T16:30:00Z"
$wsSum.Cells.Item(14, 4).Value = "2025-12-04T16:30:00Z"
$wsSum.Cells.Item(15, 4).Value = "2025-12-04T16:30:00Z"
$wsSum.Cells.Item(15, 5).Value = 31
$wsSum.Cells.Item(15, 6).Value = 964
$wsSum.Cells.Item(15, 7).Value = 1037
$wsSum.Cells.Item(15, 8).Value = 31.1
$wsSum.Cells.Item(16, 4).Value = "2025-12-04T16:30:00Z"
$wsSum.Cells.Item(17, 4).Value = "2025-12-04T16:30:00Z"
$wsSum.Cells.Item(18, 4).Value = "2025-12-04T16:30:00Z"
$wsSum.Cells.Item(19, 4).Value = "2025-12-04T16:30:00Z"
$wsSum.Cells.Item(20, 4).Value = "2025-12-04T16:30:00Z"
$wsSum.Cells.Item(21, 4).Value = "2025-12-04T16:30:00Z"
$wsSum.Cells.Item(22, 4).Value = "2025-12-04T16:30:00Z"
$wsSum.Cells.Item(22, 5).Value = 32
$wsSum.Cells.Item(22, 6).Value = 887
$wsSum.Cells.Item(22, 7).Value = 1100
$wsSum.Cells.Item(22, 8).Value = 27.7
$wsSum.Cells.Item(22, 9).Value = 34.4
$wsSum.Cells.Item(23, 4).Value = "2025-12-04T16:30:00Z"

# ---------------------------------------------------------------------
# Meta_ext: bump as_of_utc + build_version
# ---------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Meta_ext")
$wsMeta.Cells.Item(2, 2).Value = "2025-12-04T16:30:00Z"
$wsMeta.Cells.Item(2, 4).Value = 33
